$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, shifting existing rows 194..330 down to 195..331
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new data record
$ws.Cells.Item(194, 1).Value = 10
$ws.Cells.Item(194, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(194, 3).Value = "La Araucanía"
$ws.Cells.Item(194, 4).Value = 44957
$ws.Cells.Item(194, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(194, 5).Value = 9
$ws.Cells.Item(194, 6).Value = 100112039
$ws.Cells.Item(194, 7).Value = "Ciboulette"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 70
$ws.Cells.Item(194, 11).Value = 5000
$ws.Cells.Item(194, 12).Value = 6000
$ws.Cells.Item(194, 13).Value = 5429
$ws.Cells.Item(194, 14).Value = "$/docena de atados"
$ws.Cells.Item(194, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(194, 16).Value = 1810
$ws.Cells.Item(194, 17).Value = 3
$ws.Cells.Item(194, 18).Value = "Hortaliza"
